$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Subscript-3 character used in row 21 price (0.0₃0825)
$sub3 = [char]0x2083

function Set-TextValue($range, $value) {
    # Force the cell to Text format so Excel does not silently
    # reinterpret a numeric-looking string (e.g. "57.80") as a
    # floating point number and drop the trailing zero, then
    # restore the default style so no stray formatting sticks.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "37.872.88"
$ws.Range("E2").Value = "  +1.77%  "

Set-TextValue $ws.Range("D3") "2.104.86"
$ws.Range("E3").Value = "  +2.24%  "

$ws.Range("E4").Value = "  +0.03%  "

Set-TextValue $ws.Range("D5") "233.63"
$ws.Range("E5").Value = "  +0.64%  "

$ws.Range("E6").Value = "  +0.42%  "

$ws.Range("E7").Value = "  -0.04%  "

Set-TextValue $ws.Range("D8") "57.80"
$ws.Range("E8").Value = "  +1.15%  "

$ws.Range("E9").Value = "  +2.11%  "

$ws.Range("E11").Value = "  +3.45%  "

Set-TextValue $ws.Range("D12") "2.403.36"
$ws.Range("E12").Value = "  +1.86%  "

Set-TextValue $ws.Range("D13") "14.59"
$ws.Range("E13").Value = "  +0.42%  "

Set-TextValue $ws.Range("D14") "21.23"
$ws.Range("E14").Value = "  +2.11%  "

Set-TextValue $ws.Range("D15") "0.775"
$ws.Range("E15").Value = "  -0.42%  "

Set-TextValue $ws.Range("D16") "5.26"
$ws.Range("E16").Value = "  +2.35%  "

Set-TextValue $ws.Range("D17") "2.104.47"
$ws.Range("E17").Value = "  +2.36%  "

Set-TextValue $ws.Range("D18") "37.844.81"
$ws.Range("E18").Value = "  +1.87%  "

Set-TextValue $ws.Range("D19") "6.18"
$ws.Range("E19").Value = "  -3.28%  "

Set-TextValue $ws.Range("D20") "71.09"
$ws.Range("E20").Value = "  +2.75%  "

$d21val = '{0}{1}{2}' -f '0.0', $sub3, '0825'
Set-TextValue $ws.Range("D21") $d21val
$ws.Range("E21").Value = "  +2.14%  "

Set-TextValue $ws.Range("D22") "228.07"
$ws.Range("E22").Value = "  +1.12%  "

$ws.Range("E23").Value = "  +0.03%  "

Set-TextValue $ws.Range("D24") "2.41"
$ws.Range("E24").Value = "  -0.09%  "

$ws.Range("E25").Value = "  +0.06%  "

Set-TextValue $ws.Range("D26") "168.30"
$ws.Range("E26").Value = "  +1.69%  "

$ws.Range("E27").Value = "  +10.28%  "

Set-TextValue $ws.Range("D28") "8.99"
$ws.Range("E28").Value = "  +2.47%  "

$ws.Range("E29").Value = "  -1.37%  "

Set-TextValue $ws.Range("D30") "19.52"
$ws.Range("E30").Value = "  +2.65%  "

$ws.Range("E31").Value = "  +0.86%  "

$ws.Range("E32").Value = "  +4.42%  "

Set-TextValue $ws.Range("D33") "0.0632"
$ws.Range("E33").Value = "  +2.41%  "

Set-TextValue $ws.Range("D34") "4.63"
$ws.Range("E34").Value = "  +0.42%  "

Set-TextValue $ws.Range("D35") "2.55"
$ws.Range("E35").Value = "  +1.80%  "

Set-TextValue $ws.Range("D36") "3.45"
$ws.Range("E36").Value = "  +6.04%  "

$ws.Range("E37").Value = "  +4.46%  "

$ws.Range("E38").Value = "  -0.02%  "

Set-TextValue $ws.Range("D39") "5.42"
$ws.Range("E39").Value = "  -4.34%  "

Set-TextValue $ws.Range("D40") "0.0993"
$ws.Range("E40").Value = "  +6.71%  "

Set-TextValue $ws.Range("D41") "2.93"
$ws.Range("E41").Value = "  -0.51%  "

Set-TextValue $ws.Range("D42") "97.70"
$ws.Range("E42").Value = "  +1.43%  "

Set-TextValue $ws.Range("D43") "0.0215"
$ws.Range("E43").Value = "  +1.93%  "

Set-TextValue $ws.Range("D44") "1.457.17"
$ws.Range("E44").Value = "  -0.71%  "

$ws.Range("E45").Value = "  -0.25%  "

$ws.Range("E46").Value = "  +4.26%  "

Set-TextValue $ws.Range("D47") "15.71"
$ws.Range("E47").Value = "  +4.24%  "

Set-TextValue $ws.Range("D48") "4.07"
$ws.Range("E48").Value = "  -7.23%  "

Set-TextValue $ws.Range("D49") "7.39"
$ws.Range("E49").Value = "  +3.26%  "

$ws.Range("E50").Value = "  +2.38%  "

Set-TextValue $ws.Range("D51") "2.299.24"
$ws.Range("E51").Value = "  +2.31%  "
